# Generate Report for Handoff
# Refresh the localization-status report: new source-doc id
# (28c3d9fa-a100-4511-8337-8a077d2c0754 -> 2420c93c-477f-43fa-99c0-2b899dc724e6),
# new handoff-file hash (806e2b62f210b9bb460f53b9bf54eceb8e241cc5 ->
# b0811867e7215e3e5a42a28afc3d9436e40d1178), and refreshed handoff timestamps.

$wb = $excel.ActiveWorkbook

$oldId  = "28c3d9fa-a100-4511-8337-8a077d2c0754"
$newId  = "2420c93c-477f-43fa-99c0-2b899dc724e6"
$oldSha = "806e2b62f210b9bb460f53b9bf54eceb8e241cc5"
$newSha = "b0811867e7215e3e5a42a28afc3d9436e40d1178"

$mdName    = "$newId.md"
$zhXlfName = "$newId.$newSha.zh-cn.xlf"
$deXlfName = "$newId.$newSha.de-de.xlf"

# ---- Overview sheet ----
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("A2").Value = $mdName
$wsOverview.Range("D2").Value = "2016-03-22 04:57:11"
foreach ($h in $wsOverview.Hyperlinks) {
    $h.TextToDisplay = $mdName
}

# ---- zh-cn sheet ----
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("A2").Value = $mdName
$wsZhCn.Range("D2").Value = $zhXlfName
$wsZhCn.Range("E2").Value = "2016-03-22 04:57:07"
foreach ($h in $wsZhCn.Hyperlinks) {
    if ($h.TextToDisplay -like "*.md") {
        $h.TextToDisplay = $mdName
    } else {
        $h.TextToDisplay = $zhXlfName
    }
}

# ---- de-de sheet ----
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("A2").Value = $mdName
$wsDeDe.Range("D2").Value = $deXlfName
$wsDeDe.Range("E2").Value = "2016-03-22 04:57:11"
foreach ($h in $wsDeDe.Hyperlinks) {
    if ($h.TextToDisplay -like "*.md") {
        $h.TextToDisplay = $mdName
    } else {
        $h.TextToDisplay = $deXlfName
    }
}
